{"js": "// Highlight quantitative metrics (percentages, dollar amounts, large\n// numbers) in specific resume bullet paragraphs using hybrid bold + color\n// (#2C3E50) formatting, matching the author's \"quantitative metrics\n// highlighting\" commit.\n//\n// Approach: for each target paragraph (identified by its stable index in\n// the document body), locate each metric substring with Range.search()\n// (scoped to that paragraph, so matches elsewhere in the doc are not\n// touched) and apply bold + color directly to the found run. Word/Office.js\n// automatically splits the paragraph's runs around the matched text, which\n// reproduces the <w:r> run-splitting seen in the diff.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// paragraphIndex -> ordered list of metric substrings to bold+color.\nconst targets = [\n  [9, [\"23%\", \"64%\"]],\n  [11, [\"\u00b14.2%\", \"\u00b12.1%\", \"71%\", \"87%\"]],\n  [12, [\"73.5%\", \"$4.7M\"]],\n  [13, [\"$2\"]],\n  [49, [\"73.5%\"]],\n  [50, [\"$4.7M\"]],\n  [52, [\"178%\"]],\n];\n\nconst HIGHLIGHT_COLOR = \"#2C3E50\";\n\nfor (const [paraIndex, metrics] of targets) {\n  const paragraph = paragraphs.items[paraIndex];\n  for (const metric of metrics) {\n    const found = paragraph.search(metric, { matchCase: true });\n    found.load(\"items\");\n    await context.sync();\n\n    if (found.items.length === 0) continue;\n\n    const hit = found.items[0];\n    hit.font.bold = true;\n    hit.font.color = HIGHLIGHT_COLOR;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Highlight quantitative metrics (percentages, dollar amounts, large\n# numbers) in specific resume bullet paragraphs using hybrid bold + color\n# (#2C3E50) formatting, matching the author's \"quantitative metrics\n# highlighting\" commit.\n#\n# Approach: for each target paragraph (1-based Paragraphs index), use\n# Find.Execute scoped to that paragraph's own Range to locate each metric\n# substring, then bold + color the found sub-range directly. Word's COM\n# model auto-splits runs around the formatted sub-range, reproducing the\n# <w:r> run-splitting seen in the diff.\n\n$d = $word.ActiveDocument\n\n# 1-based paragraph index -> ordered list of metric substrings to bold+color.\n$targets = [ordered]@{\n    10 = @(\"23%\", \"64%\")\n    12 = @(\"\u00b14.2%\", \"\u00b12.1%\", \"71%\", \"87%\")\n    13 = @(\"73.5%\", \"$4.7M\")\n    14 = @(\"$2\")\n    50 = @(\"73.5%\")\n    51 = @(\"$4.7M\")\n    53 = @(\"178%\")\n}\n\n$highlightColor = \"2C3E50\"\n\nforeach ($paraIndex in $targets.Keys) {\n    $paragraph = $d.Paragraphs.Item($paraIndex)\n    $paraRange = $paragraph.Range\n\n    foreach ($metric in $targets[$paraIndex]) {\n        # Re-derive a fresh search range scoped to the paragraph's own text\n        # each time so consecutive finds start from the paragraph start and\n        # never drift into a neighboring paragraph.\n        $searchRange = $d.Range($paraRange.Start, $paraRange.End)\n        $find = $searchRange.Find\n        $find.ClearFormatting()\n        $find.Text = $metric\n        $find.MatchCase = $true\n        $find.MatchWholeWord = $false\n        $find.MatchWildcards = $false\n        $find.Forward = $true\n        $find.Wrap = 0\n\n        $found = $find.Execute()\n        if ($found) {\n            $searchRange.Font.Bold = $true\n            $searchRange.Font.Color = $highlightColor\n        }\n    }\n}\n"}
